# Auto-generated edit script: updates market price / profit columns (H-N)
# on multiple worksheets (ALC, ARM, CRP, CUL, LTW, WVR) per upstream data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 2043
$ws.Range("I6").Value = 2120
$ws.Range("J6").Value = 2000.2222
$ws.Range("K6").Value = 6360
$ws.Range("L6").Value = 6000.6666
$ws.Range("M6").Value = -6248
$ws.Range("N6").Value = -6224.6666
# Row 137
$ws.Range("H137").Value = 48755
$ws.Range("I137").Value = 84277.914
$ws.Range("J137").Value = 1391.1111
$ws.Range("K137").Value = 252833.742
$ws.Range("L137").Value = 4173.3333
$ws.Range("M137").Value = -250283.742
$ws.Range("N137").Value = -9273.3333
# Row 138
$ws.Range("H138").Value = 6111663.5
$ws.Range("I138").Value = 1531.174
$ws.Range("J138").Value = 8493580
$ws.Range("K138").Value = 4593.522
$ws.Range("L138").Value = 25480740
$ws.Range("M138").Value = 546.4780000000001
$ws.Range("N138").Value = -25491020
# Row 141
$ws.Range("H141").Value = 2569.9773
$ws.Range("I141").Value = 2083.4614
$ws.Range("J141").Value = 3272.7222
$ws.Range("K141").Value = 6250.3842
$ws.Range("L141").Value = 9818.1666
$ws.Range("M141").Value = -1070.3842
$ws.Range("N141").Value = -20178.1666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 10
$ws.Range("H10").Value = 41802
$ws.Range("I10").Value = 10000
$ws.Range("J10").Value = 49752.5
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 49752.5
$ws.Range("M10").Value = -9830
$ws.Range("N10").Value = -50092.5
# Row 31
$ws.Range("H31").Value = 8180.6665
$ws.Range("I31").Value = 2371
$ws.Range("J31").Value = 19800
$ws.Range("K31").Value = 2371
$ws.Range("L31").Value = 19800
$ws.Range("M31").Value = -2077
$ws.Range("N31").Value = -20388
# Row 122
$ws.Range("H122").Value = 747.85
$ws.Range("I122").Value = 591.0625
$ws.Range("J122").Value = 1375
$ws.Range("K122").Value = 1773.1875
$ws.Range("L122").Value = 4125
$ws.Range("M122").Value = 676.8125
$ws.Range("N122").Value = -9025
# Row 124
$ws.Range("H124").Value = 22000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 22000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 22000
$ws.Range("N124").Value = -31820

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1581.4259
$ws.Range("I58").Value = 638.70734
$ws.Range("J58").Value = 4554.615
$ws.Range("K58").Value = 638.70734
$ws.Range("L58").Value = 4554.615
$ws.Range("M58").Value = -435.70734
$ws.Range("N58").Value = -4960.615
# Row 74
$ws.Range("H74").Value = 12784.454
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 12784.454
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 12784.454
$ws.Range("N74").Value = -14532.454
# Row 77
$ws.Range("H77").Value = 12784.454
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 12784.454
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 38353.362
$ws.Range("N77").Value = -47089.362
# Row 132
$ws.Range("H132").Value = 1462.2916
$ws.Range("I132").Value = 843.3570999999999
$ws.Range("J132").Value = 2328.8
$ws.Range("K132").Value = 2530.0713
$ws.Range("L132").Value = 6986.400000000001
$ws.Range("M132").Value = -0.07129999999961001
$ws.Range("N132").Value = -12046.4
# Row 134
$ws.Range("H134").Value = 987.86957
$ws.Range("I134").Value = 993.70734
$ws.Range("J134").Value = 940
$ws.Range("K134").Value = 2981.12202
$ws.Range("L134").Value = 2820
$ws.Range("M134").Value = -446.1220200000002
$ws.Range("N134").Value = -7890
# Row 136
$ws.Range("H136").Value = 1581.4259
$ws.Range("I136").Value = 638.70734
$ws.Range("J136").Value = 4554.615
$ws.Range("K136").Value = 1916.12202
$ws.Range("L136").Value = 13663.845
$ws.Range("M136").Value = 633.8779799999998
$ws.Range("N136").Value = -18763.845

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 6110
$ws.Range("I5").Value = 899.1667
$ws.Range("J5").Value = 15042.857
$ws.Range("K5").Value = 2697.5001
$ws.Range("L5").Value = 45128.571
$ws.Range("M5").Value = -2585.5001
$ws.Range("N5").Value = -45352.571
# Row 7
$ws.Range("H7").Value = 300
$ws.Range("I7").Value = 233.33333
$ws.Range("J7").Value = 333.33334
$ws.Range("K7").Value = 699.99999
$ws.Range("L7").Value = 1000.00002
$ws.Range("M7").Value = -587.99999
$ws.Range("N7").Value = -1224.00002
# Row 92 (N92 is removed/replaced by new M92 in this row)
$ws.Range("H92").Value = 1000
$ws.Range("I92").Value = 1000
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 3000
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -1752
$ws.Range("N92").ClearContents()
# Row 113
$ws.Range("H113").Value = 616.6667
$ws.Range("I113").Value = 447.6
$ws.Range("J113").Value = 661.1579
$ws.Range("K113").Value = 1342.8
$ws.Range("L113").Value = 1983.4737
$ws.Range("M113").Value = 827.1999999999998
$ws.Range("N113").Value = -6323.4737
# Row 122
$ws.Range("H122").Value = 452.97296
$ws.Range("I122").Value = 220.09091
$ws.Range("J122").Value = 794.5333000000001
$ws.Range("K122").Value = 1980.81819
$ws.Range("L122").Value = 7150.7997
$ws.Range("M122").Value = 469.18181
$ws.Range("N122").Value = -12050.7997
# Row 135
$ws.Range("H135").Value = 6110
$ws.Range("I135").Value = 899.1667
$ws.Range("J135").Value = 15042.857
$ws.Range("K135").Value = 8092.5003
$ws.Range("L135").Value = 135385.713
$ws.Range("M135").Value = -5557.5003
$ws.Range("N135").Value = -140455.713
# Row 137
$ws.Range("H137").Value = 12091947
$ws.Range("I137").Value = 15152557
$ws.Range("J137").Value = 9286389
$ws.Range("K137").Value = 45457671
$ws.Range("L137").Value = 27859167
$ws.Range("M137").Value = -45452571
$ws.Range("N137").Value = -27869367

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 8427
$ws.Range("I9").Value = 254
$ws.Range("J9").Value = 35670.332
$ws.Range("K9").Value = 254
$ws.Range("L9").Value = 35670.332
$ws.Range("M9").Value = -30
$ws.Range("N9").Value = -36118.332
# Row 30
$ws.Range("H30").Value = 19133.5
$ws.Range("I30").Value = 758
$ws.Range("J30").Value = 37509
$ws.Range("K30").Value = 758
$ws.Range("L30").Value = 37509
$ws.Range("M30").Value = -650
$ws.Range("N30").Value = -37725
# Row 34
$ws.Range("H34").Value = 40000
$ws.Range("I34").Value = 40000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 40000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -39828
# Row 35
$ws.Range("H35").Value = 12519.25
$ws.Range("I35").Value = 1031
$ws.Range("J35").Value = 31666.334
$ws.Range("K35").Value = 1031
$ws.Range("L35").Value = 31666.334
$ws.Range("M35").Value = -695
$ws.Range("N35").Value = -32338.334
# Row 127
$ws.Range("H127").Value = 30000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 30000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 30000
$ws.Range("N127").Value = -39920
# Row 136
$ws.Range("H136").Value = 264794.56
$ws.Range("I136").Value = 313727.28
$ws.Range("J136").Value = 3820
$ws.Range("K136").Value = 941181.8400000001
$ws.Range("L136").Value = 11460
$ws.Range("M136").Value = -938631.8400000001
$ws.Range("N136").Value = -16560

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 80
$ws.Range("H80").Value = 47383.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 47383.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 47383.5
$ws.Range("N80").Value = -49379.5
# Row 83
$ws.Range("H83").Value = 47383.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 47383.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 142150.5
$ws.Range("N83").Value = -152134.5
